# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" sheets, which hold duplicate data, to match the newly
# generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2190
    "F4" = 1570
    "F5" = 7344
    "F7" = 180
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
